$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - LoRa Ra-01SH: new supplier link + updated price
$ws.Range("B7").Value = "https://www.aliexpress.com/item/4001148156263.html"
$ws.Range("C7").Value = 4.82

# Row 11 - USB-C: new connector part link, add price, and add digi-pote symbol link in column H
$ws.Range("B11").Value = "https://www.mouser.ee/ProductDetail/TE-Connectivity/2305018-2?qs=EU6FO9ffTwfL23TDwkN0SQ%3D%3D"
$ws.Range("C11").Value = 2.55
$ws.Range("H11").Value = "https://www.aliexpress.com/item/1005003210911840.html"

# Row 15 - RGB LED: add unit price formula (9.86/100) and quantity
$ws.Range("C15").Formula = "=9.86/100"
$ws.Range("D15").Value = 12

# Row 21 - Digi pote: USB_source voltage regulator part link update
$ws.Range("B21").Value = "https://www.mouser.ee/ProductDetail/Microchip-Technology-Atmel/MCP4019T-502E-LT?qs=%2FsslhGPpiOTvaIZO0TEW6Q%3D%3D"

# Row 25 - new component: Proximity card
$ws.Range("A25").Value = "Proximity card"
$ws.Range("B25").Value = "https://www.aliexpress.com/item/33016782798.html"

# Update view state to match the saved selection/scroll position
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("G21").Select()
